{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block (and\n// the blank paragraph that separates it from the bibliography above it),\n// which the site generator no longer emits in this build.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyrightText =\n  \"\\u00A9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\n// Locate the \"Ver no Jupiter ...\" paragraph.\nlet jupiterIdx = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === jupiterText) {\n    jupiterIdx = i;\n    break;\n  }\n}\nif (jupiterIdx === -1) {\n  throw new Error('Paragraph \"' + jupiterText + '\" not found.');\n}\n\n// The copyright notice immediately follows it.\nconst copyrightIdx = jupiterIdx + 1;\nif (\n  copyrightIdx >= paragraphs.items.length ||\n  paragraphs.items[copyrightIdx].text !== copyrightText\n) {\n  throw new Error(\"Copyright paragraph did not immediately follow the Jupiter paragraph.\");\n}\n\n// The blank separator paragraph right before \"Ver no Jupiter ...\" goes too.\nlet deleteFromIdx = jupiterIdx;\nconst precedingIdx = jupiterIdx - 1;\nif (precedingIdx >= 0 && paragraphs.items[precedingIdx].text === \"\") {\n  deleteFromIdx = precedingIdx;\n}\n\n// Delete back-to-front so earlier indices stay valid.\nfor (let i = copyrightIdx; i >= deleteFromIdx; i--) {\n  paragraphs.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block (and\n# the blank paragraph that separates it from the bibliography above it),\n# which the site generator no longer emits in this build.\n\n$d = $word.ActiveDocument\n\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\nfunction Get-ParaText($para) {\n    return $para.Range.Text.TrimEnd([char]7, [char]13, [char]10)\n}\n\n$count = $d.Paragraphs.Count\n$jupiterIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ((Get-ParaText $d.Paragraphs.Item($i)) -eq $jupiterText) {\n        $jupiterIdx = $i\n        break\n    }\n}\nif ($jupiterIdx -eq -1) {\n    throw \"Paragraph '$jupiterText' not found.\"\n}\n\n# The copyright notice immediately follows it.\n$copyrightIdx = $jupiterIdx + 1\nif ($copyrightIdx -gt $count -or (Get-ParaText $d.Paragraphs.Item($copyrightIdx)) -ne $copyrightText) {\n    throw \"Copyright paragraph did not immediately follow the Jupiter paragraph.\"\n}\n\n# The blank separator paragraph right before \"Ver no Jupiter ...\" goes too.\n$deleteFromIdx = $jupiterIdx\n$precedingIdx = $jupiterIdx - 1\nif ($precedingIdx -ge 1 -and (Get-ParaText $d.Paragraphs.Item($precedingIdx)) -eq \"\") {\n    $deleteFromIdx = $precedingIdx\n}\n\n# Delete back-to-front so earlier indices stay valid.\nfor ($i = $copyrightIdx; $i -ge $deleteFromIdx; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
